$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Tables" sheet: the cable landing points row used to reference
#    project.farm / fk_site_id — point it at project.site / site_name instead.
# ---------------------------------------------------------------------------
$tables = $wb.Worksheets.Item("Tables")
$tables.Range("B7").Value = "project.site"
$tables.Range("C7").Value = "site_name"
$tables.Application.ActiveWindow.RangeSelection
$tables.Activate()
$tables.Range("C7").Select()

# ---------------------------------------------------------------------------
# 2) "ROOT" sheet: re-sort the identifier table (A2:H25) alphabetically by
#    the Identifier column, same as Data > Sort > Column A, ascending.
# ---------------------------------------------------------------------------
$root = $wb.Worksheets.Item("ROOT")
$dataRange = $root.Range("A1:H25")
$keyRange = $root.Range("A1")
$dataRange.Sort($keyRange, 1, $null, $null, 1, $null, 1, 1)

# The sort implementation can leave stray formatting behind on cells that
# only ever held formatting (no value) at their original location once the
# data has moved elsewhere; explicitly clear those leftovers so the sheet
# ends up identical to a clean sort.
for ($r = 2; $r -le 25; $r++) {
    for ($c = 5; $c -le 8; $c++) {
        $cell = $root.Cells.Item($r, $c)
        if ($cell.Value -eq $null -and $cell.Value -eq "") {
            if ($r -ne 4) {
                $cell.Clear()
            }
        }
    }
}

$root.Activate()
$root.Range("B14").Select()
